$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05033064946979948
$ws.Range("D2").Value = 0.03515554760299722
$ws.Range("E2").Value = 0.06527041018029678
$ws.Range("F2").Value = 1.567644562272051
$ws.Range("G2").Value = 1.62744201868253
$ws.Range("H2").Value = 1.116008442887392
$ws.Range("M2").Value = 1.467739151372143
$ws.Range("N2").Value = 1.528808783980793
$ws.Range("C3").Value = 0.04456577438858744
$ws.Range("D3").Value = 0.03163531534735142
$ws.Range("E3").Value = 0.06185255765556619
$ws.Range("F3").Value = 1.43847920142619
$ws.Range("G3").Value = 1.470621640639223
$ws.Range("H3").Value = 1.050735519753488
$ws.Range("M3").Value = 1.29305177298744
$ws.Range("N3").Value = 1.414726733157437
$ws.Range("C4").Value = 0.04105264017910315
$ws.Range("D4").Value = 0.02951204739292024
$ws.Range("E4").Value = 0.05976916569155932
$ws.Range("F4").Value = 1.360514606746648
$ws.Range("G4").Value = 1.375670021507204
$ws.Range("H4").Value = 1.011663489069832
$ws.Range("M4").Value = 1.185896222750287
$ws.Range("N4").Value = 1.344947110540602
$ws.Range("C5").Value = 0.03962734940294865
$ws.Range("D5").Value = 0.02865597860972713
$ws.Range("E5").Value = 0.05892427954307777
$ws.Range("F5").Value = 1.329071281758758
$ws.Range("G5").Value = 1.337300500523241
$ws.Range("H5").Value = 0.9959886964501834
$ws.Range("M5").Value = 1.142252676382824
$ws.Range("N5").Value = 1.316581670004467
$ws.Range("C6").Value = 0.03939105435665624
$ws.Range("D6").Value = 0.0285143714120224
$ws.Range("E6").Value = 0.05878424387303127
$ws.Range("F6").Value = 1.323869657342925
$ws.Range("G6").Value = 1.330948491926819
$ws.Range("H6").Value = 0.9934006899386247
$ws.Range("M6").Value = 1.13500702298974
$ws.Range("N6").Value = 1.311875950818063
$ws.Range("C7").Value = 0.04103339297019204
$ws.Range("D7").Value = 0.02950046552759744
$ws.Range("E7").Value = 0.05975775421868335
$ws.Range("F7").Value = 1.360089236509907
$ws.Range("G7").Value = 1.375151260564195
$ws.Range("H7").Value = 1.01145109896342
$ws.Range("M7").Value = 1.185307540181185
$ws.Range("N7").Value = 1.344564275358579
$ws.Range("C8").Value = 0.04833721424105875
$ws.Range("D8").Value = 0.03393359516844185
$ws.Range("E8").Value = 0.06408895984083074
$ws.Range("F8").Value = 1.522823436754095
$ws.Range("G8").Value = 1.573085528772083
$ws.Range("H8").Value = 1.093290197834278
$ws.Range("M8").Value = 1.407483783954348
$ws.Range("N8").Value = 1.489419649574415
$ws.Range("C9").Value = 0.06288618654748745
$ws.Range("D9").Value = 0.04294929708620998
$ws.Range("E9").Value = 0.07269019997768567
$ws.Range("F9").Value = 1.853075720278781
$ws.Range("G9").Value = 1.972412569801008
$ws.Range("H9").Value = 1.262019712020276
$ws.Range("M9").Value = 1.844140479807237
$ws.Range("N9").Value = 1.7754812217986
$ws.Range("C10").Value = 0.07373636131187311
$ws.Range("D10").Value = 0.0497976118975032
$ws.Range("E10").Value = 0.07905848501506796
$ws.Range("F10").Value = 2.103183951280442
$ws.Range("G10").Value = 2.273444457487756
$ws.Range("H10").Value = 1.391393058591348
$ws.Range("M10").Value = 2.165791948797619
$ws.Range("N10").Value = 1.986735673655517
$ws.Range("C11").Value = 0.07871255083101403
$ws.Range("D11").Value = 0.05296806861487369
$ws.Range("E11").Value = 0.08196263216468935
$ws.Range("F11").Value = 2.218740759894956
$ws.Range("G11").Value = 2.412238016540527
$ws.Range("H11").Value = 1.451508328980537
$ws.Range("M11").Value = 2.312360320880202
$ws.Range("N11").Value = 2.083050618322829
$ws.Range("C12").Value = 0.08060312708707329
$ws.Range("D12").Value = 0.05417707685026585
$ws.Range("E12").Value = 0.08306307578183691
$ws.Range("F12").Value = 2.262767672201988
$ws.Range("G12").Value = 2.465076962121202
$ws.Range("H12").Value = 1.474460930874216
$ws.Range("M12").Value = 2.367901839402748
$ws.Range("N12").Value = 2.119550861324512
$ws.Range("C13").Value = 0.08019567549469286
$ws.Range("D13").Value = 0.05391631252516049
$ws.Range("E13").Value = 0.08282604938253968
$ws.Range("F13").Value = 2.253273574270366
$ws.Range("G13").Value = 2.453684425982658
$ws.Range("H13").Value = 1.46950919653176
$ws.Range("M13").Value = 2.355938173348136
$ws.Range("N13").Value = 2.111688687718129
$ws.Range("C14").Value = 0.0788679631061342
$ws.Range("D14").Value = 0.05306736270148349
$ws.Range("E14").Value = 0.08205315405209035
$ws.Range("F14").Value = 2.222357438962518
$ws.Range("G14").Value = 2.416579396090413
$ws.Range("H14").Value = 1.453392841597179
$ws.Range("M14").Value = 2.316928946615036
$ws.Range("N14").Value = 2.086052973438768
$ws.Range("C15").Value = 0.07805552036218444
$ws.Range("D15").Value = 0.05254846861163287
$ws.Range("E15").Value = 0.08157981556399818
$ws.Range("F15").Value = 2.203455680688762
$ws.Range("G15").Value = 2.393888501428535
$ws.Range("H15").Value = 1.443545826211221
$ws.Range("M15").Value = 2.293039893216729
$ws.Range("N15").Value = 2.070353900862926
$ws.Range("C16").Value = 0.0734120017692419
$ws.Range("D16").Value = 0.04959156665674413
$ws.Range("E16").Value = 0.07886880882583114
$ws.Range("F16").Value = 2.095668897786908
$ws.Range("G16").Value = 2.264412477397684
$ws.Range("H16").Value = 1.387490360021104
$ws.Range("M16").Value = 2.156218650674901
$ws.Range("N16").Value = 1.980445323152196
$ws.Range("C17").Value = 0.07057397513804631
$ws.Range("D17").Value = 0.04779207247688078
$ws.Range("E17").Value = 0.07720730190693104
$ws.Range("F17").Value = 2.030009814933891
$ws.Range("G17").Value = 2.185467793540454
$ws.Range("H17").Value = 1.35343025182425
$ws.Range("M17").Value = 2.072349193459758
$ws.Range("N17").Value = 1.925342075430422
$ws.Range("C18").Value = 0.0689453949251515
$ws.Range("D18").Value = 0.04676220902967998
$ws.Range("E18").Value = 0.0762523427935804
$ws.Range("F18").Value = 1.992411383770332
$ws.Range("G18").Value = 2.140234395058826
$ws.Range("H18").Value = 1.333958218775535
$ws.Range("M18").Value = 2.024132766020216
$ws.Range("N18").Value = 1.893668480661944
$ws.Range("C19").Value = 0.06839462264724716
$ws.Range("D19").Value = 0.04641438586747881
$ws.Range("E19").Value = 0.07592913944989022
$ws.Range("F19").Value = 1.979709514567702
$ws.Range("G19").Value = 2.124948507661827
$ws.Range("H19").Value = 1.327385464356098
$ws.Range("M19").Value = 2.007811365063844
$ws.Range("N19").Value = 1.882947929325383
$ws.Range("C20").Value = 0.07087569412989581
$ws.Range("D20").Value = 0.04798309429806125
$ws.Range("E20").Value = 0.07738410250898653
$ws.Range("F20").Value = 2.036981973777586
$ws.Range("G20").Value = 2.193853530353806
$ws.Range("H20").Value = 1.357043702508236
$ws.Range("M20").Value = 2.08127482712186
$ws.Range("N20").Value = 1.931205827090196
$ws.Range("C21").Value = 0.07925777293445435
$ws.Range("D21").Value = 0.0533164871585825
$ws.Range("E21").Value = 0.08228015559786428
$ws.Range("F21").Value = 2.231430885989596
$ws.Range("G21").Value = 2.427470304236124
$ws.Range("H21").Value = 1.45812144246463
$ws.Range("M21").Value = 2.328385806561442
$ws.Range("N21").Value = 2.09358207100297
$ws.Range("C22").Value = 0.08477228150611893
$ws.Range("D22").Value = 0.05685151546894929
$ws.Range("E22").Value = 0.08548397173694156
$ws.Range("F22").Value = 2.360082668991964
$ws.Range("G22").Value = 2.58179625400328
$ws.Range("H22").Value = 1.525282098093953
$ws.Range("M22").Value = 2.490117769034867
$ws.Range("N22").Value = 2.199865432404977
$ws.Range("C23").Value = 0.08182562726655362
$ws.Range("D23").Value = 0.0549601212107973
$ws.Range("E23").Value = 0.08377378079791242
$ws.Range("F23").Value = 2.291271291732301
$ws.Range("G23").Value = 2.49927431499691
$ws.Range("H23").Value = 1.489334243490703
$ws.Range("M23").Value = 2.403776013123093
$ws.Range("N23").Value = 2.143126214273479
$ws.Range("C24").Value = 0.07073927753728526
$ws.Range("D24").Value = 0.04789671876309853
$ws.Range("E24").Value = 0.07730417011190127
$ws.Range("F24").Value = 2.033829394911493
$ws.Range("G24").Value = 2.190061863416759
$ws.Range("H24").Value = 1.355409720650016
$ws.Range("M24").Value = 2.077239544553095
$ws.Range("N24").Value = 1.928554806062948
$ws.Range("C25").Value = 0.05892372672171575
$ws.Range("D25").Value = 0.04047294831426029
$ws.Range("E25").Value = 0.07035352172209031
$ws.Range("F25").Value = 1.762472090930004
$ws.Range("G25").Value = 1.863102442482898
$ws.Range("H25").Value = 1.215452414927768
$ws.Range("M25").Value = 1.725884347862831
$ws.Range("N25").Value = 1.697894012675903
